$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 for the new component "C22" (alphabetically
# between C21 and D1), shifting all subsequent rows down by one.
$ws.Rows("23:23").Insert()

# Expand the "CPL" table (query table) to include the new row.
$lo = $ws.ListObjects.Item("CPL")
$lo.Resize($ws.Range("A1:E61"))

# Update the ExternalData_1 defined name that tracks the query table range.
$nm = $wb.Names.Item("ExternalData_1")
$nm.RefersTo = "='CPL'!`$A`$1:`$E`$61"

# Populate the new row's values.
$ws.Range("A23").Value = "C22"
$ws.Range("D23").Value = "top"
$ws.Range("E23").Value = 180

# Mid X / Mid Y are stored as text (the query converts them to text), so
# write them as text-producing formulas and then flatten to static text
# values via copy / paste-special so they keep the table's default style.
$ws.Range("B23").Formula = '="106.7"'
$ws.Range("C23").Formula = '="-59.5"'
$ws.Range("B23:C23").Copy()
$ws.Range("B23:C23").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Match the author's reported selection after the edit.
[void]$ws.Range("J13").Select()
